$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (everything below shifts down by one row;
# old row 52 becomes row 53, dimension grows to A1:R53 automatically).
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits the header row's bold/border style from
# Excel's default "format from above" insert behaviour; reset it back to the
# plain style used by every other data row before filling in values.
$ws.Range("A2:R2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 45092
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 100112043
$ws.Range("G2").Value = "Pepino dulce"
$ws.Range("H2").Value = "Cultivar IV Región"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("N2").Value = "$/bandeja 18 kilos"
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 778
$ws.Range("Q2").Value = 18
$ws.Range("R2").Value = "Hortaliza"
